$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2; D = "27.784.00"; E = "  +0.37%  " },
    @{ Row = 3; D = "1.865.38"; E = "  -0.47%  " },
    @{ Row = 4; D = "1.013"; E = "  +0.84%  " },
    @{ Row = 5; D = "333.61"; E = "  +0.50%  " },
    @{ Row = 6; D = "1.011"; E = "  +0.77%  " },
    @{ Row = 7; D = "0.4714"; E = "  -0.11%  " },
    @{ Row = 8; D = "0.3896"; E = "  -1.18%  " },
    @{ Row = 9; D = "46.64"; E = "  -2.71%  " },
    @{ Row = 10; D = "0.07958"; E = "  -0.88%  " },
    @{ Row = 11; D = "1.003"; E = "  -2.25%  " },
    @{ Row = 12; D = "21.53"; E = "  -2.07%  " },
    @{ Row = 13; D = "1.870.55"; E = "  +0.97%  " },
    @{ Row = 14; D = "5.972"; E = "  +0.16%  " },
    @{ Row = 15; D = "7.131"; E = "  +0.00%  " },
    @{ Row = 16; D = "1.014"; E = "  +0.68%  " },
    @{ Row = 17; D = "88.09"; E = "  +1.18%  " },
    @{ Row = 18; D = $null; E = "  +0.53%  " },
    @{ Row = 19; D = "0.00001040"; E = "  -0.74%  " },
    @{ Row = 20; D = "16.89"; E = "  -1.60%  " },
    @{ Row = 21; D = "1.010"; E = "  +0.62%  " },
    @{ Row = 22; D = "27.775.03"; E = "  +0.32%  " },
    @{ Row = 23; D = "5.448"; E = "  -1.23%  " },
    @{ Row = 24; D = "10.89"; E = "  -0.85%  " },
    @{ Row = 25; D = "2.332"; E = "  +1.06%  " },
    @{ Row = 26; D = "2.093.49"; E = "  +0.81%  " },
    @{ Row = 27; D = "157.91"; E = "  -0.16%  " },
    @{ Row = 28; D = "19.65"; E = "  -2.64%  " },
    @{ Row = 29; D = "2.083"; E = "  -0.86%  " },
    @{ Row = 30; D = "5.398"; E = "  -3.38%  " },
    @{ Row = 31; D = "120.95"; E = "  -1.03%  " },
    @{ Row = 32; D = "0.9631"; E = "  -1.20%  " },
    @{ Row = 33; D = "0.09445"; E = "  -1.03%  " },
    @{ Row = 34; D = "3.638"; E = "  +1.29%  " },
    @{ Row = 35; D = "5.287"; E = "  -0.84%  " },
    @{ Row = 36; D = "1.342"; E = "  -7.08%  " },
    @{ Row = 37; D = "0.06028"; E = "  -1.13%  " },
    @{ Row = 38; D = "0.02216"; E = "  -1.52%  " },
    @{ Row = 39; D = $null; E = "  -2.09%  " },
    @{ Row = 40; D = "8.118"; E = "  -1.39%  " },
    @{ Row = 41; D = "1.009"; E = "  +0.63%  " },
    @{ Row = 42; D = "0.5897"; E = "  -2.03%  " },
    @{ Row = 43; D = "0.1884"; E = "  -1.21%  " },
    @{ Row = 44; D = "10.22"; E = "  -0.11%  " },
    @{ Row = 45; D = "1.256"; E = "  +0.00%  " },
    @{ Row = 46; D = "0.5606"; E = "  -1.47%  " },
    @{ Row = 47; D = "11.97"; E = "  -2.29%  " },
    @{ Row = 48; D = "1.908"; E = "  -1.75%  " },
    @{ Row = 49; D = "3.304"; E = "  -2.43%  " },
    @{ Row = 50; D = "0.06776"; E = "  -1.63%  " },
    @{ Row = 51; D = "111.84"; E = "  -3.12%  " }
)

foreach ($item in $data) {
    $r = $item.Row
    if ($item.D -ne $null) {
        $ws.Range("D$r").NumberFormat = "@"
        $ws.Range("D$r").Value = $item.D
        $ws.Range("D$r").Style = "Normal"
    }
    $ws.Range("E$r").NumberFormat = "@"
    $ws.Range("E$r").Value = $item.E
    $ws.Range("E$r").Style = "Normal"
}
